$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.666.95"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "3.324.96"
$ws.Range("E3").Value = "  +5.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.15%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.320.69"
$ws.Range("E8").Value = "  +5.48%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "3.872.76"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("D17").Value = "3.321.87"
$ws.Range("E17").Value = "  +5.55%  "
$ws.Range("D18").Value = "63.757.29"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.20%  "
$ws.Range("E24").Value = "  +5.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.34%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.54%  "
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0746"
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0401"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "434.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.099.32"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.87%  "
